$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple text/string cell updates
$ws.Range("D2").Value = "66.034.57"
$ws.Range("E2").Value = "  -6.00%  "
$ws.Range("D3").Value = "3.166.76"
$ws.Range("E3").Value = "  -9.82%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  -7.10%  "
$ws.Range("E6").Value = "  -15.20%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "3.156.05"
$ws.Range("E8").Value = "  -9.98%  "
$ws.Range("E9").Value = "  -12.03%  "
$ws.Range("E10").Value = "  -14.35%  "
$ws.Range("E11").Value = "  -12.50%  "
$ws.Range("E12").Value = "  -16.85%  "
$ws.Range("E13").Value = "  -18.54%  "
$ws.Range("E14").Value = "  -14.11%  "
$ws.Range("D15").Value = "3.682.33"
$ws.Range("E15").Value = "  -9.73%  "
$ws.Range("D16").Value = "66.065.37"
$ws.Range("E16").Value = "  -5.98%  "
$ws.Range("D17").Value = "3.181.01"
$ws.Range("E18").Value = "  -6.67%  "
$ws.Range("E19").Value = "  -14.61%  "
$ws.Range("E20").Value = "  -16.91%  "
$ws.Range("E21").Value = "  -16.62%  "
$ws.Range("E22").Value = "  -15.42%  "
$ws.Range("E23").Value = "  -15.05%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("E24").Value = "  -15.07%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("E25").Value = "  -14.99%  "
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("E27").Value = "  -17.99%  "
$ws.Range("E28").Value = "  -17.41%  "
$ws.Range("E29").Value = "  -14.47%  "
$ws.Range("E30").Value = "  -15.24%  "
$ws.Range("E31").Value = "  -15.34%  "
$ws.Range("E32").Value = "  -15.25%  "
$ws.Range("E33").Value = "  -21.33%  "
$ws.Range("E34").Value = "  -16.86%  "
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("E36").Value = "  -18.30%  "
$ws.Range("E37").Value = "  -6.62%  "
$ws.Range("E38").Value = "  -16.25%  "
$ws.Range("E39").Value = "  -16.50%  "
$ws.Range("E40").Value = "  -18.76%  "
$ws.Range("E41").Value = "  -14.43%  "
$ws.Range("D42").Value = "2.842.86"
$ws.Range("E42").Value = "  -15.42%  "
$ws.Range("E43").Value = "  -28.54%  "
$ws.Range("E44").Value = "  -17.71%  "
$ws.Range("D45").Value = "0.0₃0564"
$ws.Range("E45").Value = "  -23.44%  "
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("E47").Value = "  -20.93%  "
$ws.Range("E48").Value = "  -20.37%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E49").Value = "  -8.01%  "
$ws.Range("E50").Value = "  -14.26%  "
$ws.Range("B51").Value = "Fetch.AI"
$ws.Range("C51").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("E51").Value = "  -19.98%  "

# Numeric-looking price updates forced to remain text (matches source formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.167"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.26"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.489"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.73"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000236"
$ws.Range("D14").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "522.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.745"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.53"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "517.38"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0840"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0401"
$ws.Range("D40").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.255"
$ws.Range("D44").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.04"
$ws.Range("D51").Style = "Normal"
